# Update "想去人数" (want-to-go count) figures that changed between scrapes.
# Both "展览" (sheet 1) and "全部类型" (sheet 4) list the same events, so the
# same F-column increments apply to both sheets (rows differ by 1 because
# "全部类型" also contains a 演出/performance row that 展览 does not).

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)   # 展览
$sheet4 = $wb.Worksheets.Item(4)   # 全部类型

# Updates for the "展览" sheet (row -> new value)
$sheet1.Range("F3").Value = 791
$sheet1.Range("F6").Value = 83
$sheet1.Range("F7").Value = 279
$sheet1.Range("F8").Value = 3978
$sheet1.Range("F10").Value = 4676
$sheet1.Range("F11").Value = 517
$sheet1.Range("F12").Value = 1180
$sheet1.Range("F13").Value = 78

# Same events, shifted one row down, on the "全部类型" sheet
$sheet4.Range("F3").Value = 791
$sheet4.Range("F6").Value = 83
$sheet4.Range("F8").Value = 279
$sheet4.Range("F9").Value = 3978
$sheet4.Range("F11").Value = 4676
$sheet4.Range("F12").Value = 517
$sheet4.Range("F13").Value = 1180
$sheet4.Range("F14").Value = 78
